# Invitee import/export: add "google+_link" and "linkedin_link" columns
# (with sample data) right after the existing "facebook_link" column,
# pushing the old trailing "password" column from M to O.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New header cells (row 1) ---
$ws.Range("M1").Value = "google+_link"
$ws.Range("N1").Value = "linkedin_link"
$ws.Range("O1").Value = "password"

# --- New data cells (row 2) ---
$ws.Range("M2").Value = "Coolshiv@gmail.com"
$ws.Range("N2").Value = "Coolshiv@gmail.com"
$ws.Range("O2").Value = "password"

$null = $ws.Range("N2").Select()
